$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "42.012.21"
$ws.Cells.Item(2, 5).Value = "  -1.35%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.302.67"
$ws.Cells.Item(3, 5).Value = "  -2.36%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "318.21"
$ws.Cells.Item(5, 5).Value = "  -0.92%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "104.35"
$ws.Cells.Item(6, 5).Value = "  -2.69%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.67%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.07%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.609"
$ws.Cells.Item(9, 5).Value = "  -1.86%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "39.55"
$ws.Cells.Item(10, 5).Value = "  -4.29%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0912"
$ws.Cells.Item(11, 5).Value = "  -1.37%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "8.40"
$ws.Cells.Item(12, 5).Value = "  -0.79%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +0.15%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "0.977"
$ws.Cells.Item(14, 5).Value = "  -1.94%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "15.42"
$ws.Cells.Item(15, 5).Value = "  -3.54%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "2.650.67"
$ws.Cells.Item(16, 5).Value = "  -2.38%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "2.318.83"
$ws.Cells.Item(17, 5).Value = "  -0.96%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "42.125.88"
$ws.Cells.Item(18, 5).Value = "  -0.92%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "7.68"
$ws.Cells.Item(19, 5).Value = "  +0.53%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -0.44%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "290.62"
$ws.Cells.Item(21, 5).Value = "  +13.03%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "Litecoin"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "73.95"
$ws.Cells.Item(22, 5).Value = "  -3.12%  "

# Row 23
$ws.Cells.Item(23, 2).Value = "PancakeSwap"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "3.60"
$ws.Cells.Item(23, 5).Value = "  -0.69%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -1.03%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.99"
$ws.Cells.Item(25, 5).Value = "  +6.74%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "1.00"
$ws.Cells.Item(26, 5).Value = "  -0.03%  "

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "10.96"
$ws.Cells.Item(27, 5).Value = "  -4.13%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "23.51"
$ws.Cells.Item(28, 5).Value = "  +2.49%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +1.68%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "Monero"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "164.36"
$ws.Cells.Item(30, 5).Value = "  -6.58%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "35.34"
$ws.Cells.Item(31, 5).Value = "  -3.32%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.0885"
$ws.Cells.Item(32, 5).Value = "  -0.68%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +1.65%  "

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.90"
$ws.Cells.Item(34, 5).Value = "  -2.39%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.133"
$ws.Cells.Item(35, 5).Value = "  +0.95%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -7.99%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "4.64"
$ws.Cells.Item(37, 5).Value = "  +0.35%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.0352"
$ws.Cells.Item(38, 5).Value = "  -2.91%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.86"
$ws.Cells.Item(39, 5).Value = "  +6.59%  "

# Row 40
$ws.Cells.Item(40, 5).Value = "  -6.24%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "103.28"
$ws.Cells.Item(41, 5).Value = "  +22.36%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "1.48"
$ws.Cells.Item(42, 5).Value = "  +0.77%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "70.85"
$ws.Cells.Item(43, 5).Value = "  -0.44%  "

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "0.227"
$ws.Cells.Item(44, 5).Value = "  -5.18%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  +0.35%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "117.79"
$ws.Cells.Item(46, 5).Value = "  +4.47%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "12.09"
$ws.Cells.Item(47, 5).Value = "  +0.65%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "77.79"
$ws.Cells.Item(48, 5).Value = "  +5.44%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "9.12"
$ws.Cells.Item(49, 5).Value = "  -0.09%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "5.35"
$ws.Cells.Item(50, 5).Value = "  -2.41%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  -0.08%  "
